$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# Update the dist_wos_mu (AE) and dist_wos_sigma (AF) values for rows 2-6
$ws.Range("AE2").Value = 9.453
$ws.Range("AF2").Value = 0.046

$ws.Range("AE3").Value = 9.249
$ws.Range("AF3").Value = 0.061

$ws.Range("AE4").Value = 8.243
$ws.Range("AF4").Value = 0.105

$ws.Range("AE5").Value = 8.66
$ws.Range("AF5").Value = 0.211

$ws.Range("AE6").Value = 9.412
$ws.Range("AF6").Value = 0.105

# Move the view/selection to match the author's final cursor position
$ws.Range("AG10").Select()
